$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44411
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 210
$ws.Range("N2").Value = 8000
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 8000
$ws.Range("Q2").Value = '$/bandeja 8 kilos'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 8

$ws.Range("D3").Value = 44511
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("Q3").Value = '$/caja 15 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 1467
$ws.Range("T3").Value = 15

$ws.Range("D4").Value = 44418
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("Q4").Value = '$/caja 15 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 533
$ws.Range("T4").Value = 15

$ws.Range("D5").Value = 44427
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("Q5").Value = '$/caja 15 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 467
$ws.Range("T5").Value = 15

$ws.Range("D6").Value = 44495
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 24000
$ws.Range("Q6").Value = '$/bandeja 10 kilos'
$ws.Range("R6").Value = 'China'
$ws.Range("S6").Value = 2400
$ws.Range("T6").Value = 10

$ws.Range("D7").Value = 44208
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 24000
$ws.Range("Q7").Value = '$/caja 15 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1600
$ws.Range("T7").Value = 15

$ws.Range("D8").Value = 44217
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 55
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44392
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 500
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 7000
$ws.Range("Q9").Value = '$/bandeja 8 kilos'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 875
$ws.Range("T9").Value = 8

$ws.Range("D10").Value = 44601
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 28000
$ws.Range("O10").Value = 28000
$ws.Range("P10").Value = 28000
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1556
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 44264
$ws.Range("L11").Value = 'Calibre 100'
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = '$/caja 18 kilos embalada'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1111
$ws.Range("T11").Value = 18

